$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A4: bus stop location changed from "3500 blk Kingsway ave" to "300 Block Main St"
$ws.Range("A4").Value = 'Victim was waiting at a bus stop in the 300 Block Main St when an unknown male approached her from behind and groped her buttocks. The male then followed victim in to the 3300 blk Padora Ave and fled the area. Written statement and video of suspect obtained. '

# Update A11: "On April 21, 2021" -> "Yesterday", and remove last sentence about safety plan
$ws.Range("A11").Value = 'Victim Grade 11 at Killarney Secondary School was walking to a grocery store at Champlain Square near the intersection of Pandora Ave/Kerr Street, Vancouver, BC, when all of a sudden her crotch was grabbed over her shorts from behind by an unknown male, who then fled the area on foot. Yesterday after she finished her afternoon class, victim approached her SLO in the school hallway and advised him of the circumstances. Victim did not observe the male''s face, but observed him to be wearing a black hoody with the hood up and black or grey sweatpants. Video canvass to be conducted and victim''s unwashed shorts to be obtained and submitted for forensic processing.'

# Update selection to A11
$ws.Range("A11").Select()
